# Update "想去人数" (number of people interested) values in the
# "展览" and "全部类型" worksheets to reflect the latest scrape.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F
$updates = @{
    2  = 132
    3  = 411
    4  = 11973
    5  = 1258
    11 = 341
    13 = 60
    16 = 349
    17 = 1710
    19 = 921
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
